$d = $word.ActiveDocument

# Locate the paragraph that begins "In May 1862, the Qing forces laid
# siege..." without depending on a hard-coded paragraph index.
$count = $d.Paragraphs.Count
$i = 1
$paraIndex = -1
while ($i -le $count) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.StartsWith("In May 1862, the Qing forces laid siege")) {
        $paraIndex = $i
    }
    $i = $i + 1
}

$p = $d.Paragraphs.Item($paraIndex)

# ---------------------------------------------------------------------
# 1) Add the new "5-Guangdong Highlands- " heading run at the very start
#    of the paragraph (kept as its own run).
# ---------------------------------------------------------------------
$p.Range.InsertBefore("5-Guangdong Highlands- ")

# ---------------------------------------------------------------------
# 2) Rework the paragraph's closing sentences describing Hong's death
#    and the fall of the rebellion.
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item($paraIndex)
$find = $p.Range.Duplicate()
$null = $find.Find.Execute(" died ")
$anchorStart = $find.Start

$paraContentEnd = $p.Range.End - 1
$oldTail = $d.Range($anchorStart, $paraContentEnd)
$oldTail.Delete()

function Append-Run {
    param([string]$Text)
    $pEnd = $d.Paragraphs.Item($paraIndex).Range.End - 1
    $ip = $d.Range($pEnd, $pEnd)
    $ip.InsertAfter($Text)
}

Append-Run(" died prior to the taking of the ")
Append-Run("city, though ")
Append-Run("to what")
Append-Run(" end is still a matter of dispute. He died of poison though it is unclear if he was murdered or if he died by his own hand.")
Append-Run(" His son-also named Hong-was named as his successor but was quickly caught and executed.")
Append-Run(" Regardless, the death of the Heavenly King marked the eventual downfall of the ")
Append-Run("Taiping rebellion and the continuation of Qing rule until the fall of Dynastic China in 1911. With Hong dead and ")
Append-Run("the ")
Append-Run("Heavenly Capital no longer under Taiping control, they were forced into retreat all the way back to the highlands in Guangdong, not far from where the rebellion began fifteen years prior. The Taiping rebellion was one the deadliest conflicts in human history with loss of life between 20 and 70 million. ")
